$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the data values in the TIMER2 table (row 12)
$ws.Range("B12").Value = 48
$ws.Range("D12").Value = 16
$ws.Range("E12").Value = 16

# Recalculate formulas so dependent cells (C12, B14, C14) update
$excel.Calculate()

# Update the view selection/scroll position to match the final state
$ws.Range("D12").Select()
$excel.ActiveWindow.ScrollRow = 7
